$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2444.7368
$ws.Cells.Item(32, 10).Value = 2911.6
$ws.Cells.Item(32, 12).Value = 2911.6
$ws.Cells.Item(32, 14).Value = -3563.6
$ws.Cells.Item(33, 8).Value = 311.1111
$ws.Cells.Item(33, 9).Value = 319.44
$ws.Cells.Item(33, 11).Value = 319.44
$ws.Cells.Item(33, 13).Value = -90.44
$ws.Cells.Item(51, 8).Value = 7497.15
$ws.Cells.Item(51, 9).Value = 8649.25
$ws.Cells.Item(51, 10).Value = 6729.0835
$ws.Cells.Item(51, 11).Value = 8649.25
$ws.Cells.Item(51, 12).Value = 6729.0835
$ws.Cells.Item(51, 13).Value = -8165.25
$ws.Cells.Item(51, 14).Value = -7697.0835
$ws.Cells.Item(125, 8).Value = 1097.7273
$ws.Cells.Item(125, 9).Value = 1645
$ws.Cells.Item(125, 10).Value = 976.1111
$ws.Cells.Item(125, 11).Value = 14805
$ws.Cells.Item(125, 12).Value = 8784.999899999999
$ws.Cells.Item(125, 13).Value = -12345
$ws.Cells.Item(125, 14).Value = -13704.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 26333.334
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 26333.334
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 26333.334
$ws.Cells.Item(10, 13).ClearContents()
$ws.Cells.Item(10, 14).Value = -26673.334
$ws.Cells.Item(32, 8).Value = 1298.5676
$ws.Cells.Item(32, 9).Value = 1195.75
$ws.Cells.Item(32, 11).Value = 1195.75
$ws.Cells.Item(32, 13).Value = -908.75
$ws.Cells.Item(45, 8).Value = 1504.4445
$ws.Cells.Item(45, 9).Value = 1615.75
$ws.Cells.Item(45, 10).Value = 614
$ws.Cells.Item(45, 11).Value = 1615.75
$ws.Cells.Item(45, 12).Value = 614
$ws.Cells.Item(45, 13).Value = -1238.75
$ws.Cells.Item(45, 14).Value = -1368
$ws.Cells.Item(61, 8).Value = 4660.086
$ws.Cells.Item(61, 9).Value = 3681.2222
$ws.Cells.Item(61, 10).Value = 7963.75
$ws.Cells.Item(61, 11).Value = 3681.2222
$ws.Cells.Item(61, 12).Value = 7963.75
$ws.Cells.Item(61, 13).Value = -3469.2222
$ws.Cells.Item(61, 14).Value = -8387.75
$ws.Cells.Item(110, 8).Value = 7524.875
$ws.Cells.Item(110, 9).Value = 8349.857
$ws.Cells.Item(110, 11).Value = 8349.857
$ws.Cells.Item(110, 13).Value = -6304.857
$ws.Cells.Item(122, 8).Value = 1753.4
$ws.Cells.Item(122, 9).Value = 1753.4
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 5260.200000000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -2810.200000000001
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 84900
$ws.Cells.Item(131, 10).Value = 84900
$ws.Cells.Item(131, 12).Value = 84900
$ws.Cells.Item(131, 14).Value = -94980
$ws.Cells.Item(132, 8).Value = 2600.6
$ws.Cells.Item(132, 9).Value = 2600.6
$ws.Cells.Item(132, 11).Value = 7801.799999999999
$ws.Cells.Item(132, 13).Value = -5271.799999999999
$ws.Cells.Item(136, 8).Value = 4660.086
$ws.Cells.Item(136, 9).Value = 3681.2222
$ws.Cells.Item(136, 10).Value = 7963.75
$ws.Cells.Item(136, 11).Value = 11043.6666
$ws.Cells.Item(136, 12).Value = 23891.25
$ws.Cells.Item(136, 13).Value = -8493.6666
$ws.Cells.Item(136, 14).Value = -28991.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 2187.5
$ws.Cells.Item(11, 9).Value = 2600
$ws.Cells.Item(11, 11).Value = 2600
$ws.Cells.Item(11, 13).Value = -2460
$ws.Cells.Item(20, 8).Value = 1856
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).ClearContents()
$ws.Cells.Item(35, 8).Value = 61429.2
$ws.Cells.Item(35, 9).Value = 58500
$ws.Cells.Item(35, 10).Value = 63382
$ws.Cells.Item(35, 11).Value = 58500
$ws.Cells.Item(35, 12).Value = 63382
$ws.Cells.Item(35, 13).Value = -58190
$ws.Cells.Item(35, 14).Value = -64002
$ws.Cells.Item(86, 8).Value = 5787.8887
$ws.Cells.Item(86, 10).Value = 16332.667
$ws.Cells.Item(86, 12).Value = 16332.667
$ws.Cells.Item(86, 14).Value = -18578.667
$ws.Cells.Item(89, 8).Value = 5787.8887
$ws.Cells.Item(89, 10).Value = 16332.667
$ws.Cells.Item(89, 12).Value = 81663.33499999999
$ws.Cells.Item(89, 14).Value = -92895.33499999999
$ws.Cells.Item(99, 8).Value = 968.5
$ws.Cells.Item(99, 9).Value = 742.7778
$ws.Cells.Item(99, 10).Value = 3000
$ws.Cells.Item(99, 11).Value = 742.7778
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 13).Value = 755.2222
$ws.Cells.Item(99, 14).Value = -5996
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 442.88235
$ws.Cells.Item(7, 9).Value = 355.92307
$ws.Cells.Item(7, 10).Value = 725.5
$ws.Cells.Item(7, 11).Value = 355.92307
$ws.Cells.Item(7, 12).Value = 725.5
$ws.Cells.Item(7, 13).Value = -242.92307
$ws.Cells.Item(7, 14).Value = -951.5
$ws.Cells.Item(16, 8).Value = 3550
$ws.Cells.Item(16, 9).Value = 2100
$ws.Cells.Item(16, 10).Value = 5000
$ws.Cells.Item(16, 11).Value = 2100
$ws.Cells.Item(16, 12).Value = 5000
$ws.Cells.Item(16, 13).Value = -1813
$ws.Cells.Item(16, 14).Value = -5574
$ws.Cells.Item(99, 8).Value = 2791.6875
$ws.Cells.Item(99, 9).Value = 2690.5715
$ws.Cells.Item(99, 11).Value = 2690.5715
$ws.Cells.Item(99, 13).Value = -1192.5715
$ws.Cells.Item(107, 8).Value = 698.8570999999999
$ws.Cells.Item(107, 9).Value = 718.3
$ws.Cells.Item(107, 10).Value = 650.25
$ws.Cells.Item(107, 11).Value = 718.3
$ws.Cells.Item(107, 12).Value = 650.25
$ws.Cells.Item(107, 13).Value = 1201.7
$ws.Cells.Item(107, 14).Value = -4490.25
$ws.Cells.Item(113, 8).Value = 3550
$ws.Cells.Item(113, 9).Value = 2100
$ws.Cells.Item(113, 10).Value = 5000
$ws.Cells.Item(113, 11).Value = 2100
$ws.Cells.Item(113, 12).Value = 5000
$ws.Cells.Item(113, 13).Value = 70
$ws.Cells.Item(113, 14).Value = -9340
$ws.Cells.Item(122, 8).Value = 6367
$ws.Cells.Item(122, 9).Value = 6296.6665
$ws.Cells.Item(122, 11).Value = 18889.9995
$ws.Cells.Item(122, 13).Value = -16439.9995
$ws.Cells.Item(126, 8).Value = 2791.6875
$ws.Cells.Item(126, 9).Value = 2690.5715
$ws.Cells.Item(126, 11).Value = 8071.7145
$ws.Cells.Item(126, 13).Value = -5601.7145
$ws.Cells.Item(141, 8).Value = 100163
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 100163
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 100163
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).Value = -110523
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 153.5
$ws.Cells.Item(15, 9).Value = 60
$ws.Cells.Item(15, 11).Value = 180
$ws.Cells.Item(15, 13).Value = -40
$ws.Cells.Item(29, 8).Value = 601
$ws.Cells.Item(29, 10).Value = 601
$ws.Cells.Item(29, 12).Value = 1803
$ws.Cells.Item(29, 14).Value = -2357
$ws.Cells.Item(124, 8).Value = 9500
$ws.Cells.Item(124, 9).Value = 9500
$ws.Cells.Item(124, 11).Value = 28500
$ws.Cells.Item(124, 13).Value = -23590
$ws.Cells.Item(131, 8).Value = 975.8333
$ws.Cells.Item(131, 9).Value = 882.7273
$ws.Cells.Item(131, 11).Value = 2648.1819
$ws.Cells.Item(131, 13).Value = 2391.8181
$ws.Cells.Item(132, 8).Value = 1379.1
$ws.Cells.Item(132, 9).Value = 1630.6666
$ws.Cells.Item(132, 11).Value = 14675.9994
$ws.Cells.Item(132, 13).Value = -12145.9994
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3600.3333
$ws.Cells.Item(122, 9).Value = 1981.5
$ws.Cells.Item(122, 11).Value = 5944.5
$ws.Cells.Item(122, 13).Value = -3494.5
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4687.6875
$ws.Cells.Item(7, 9).Value = 3873.2727
$ws.Cells.Item(7, 10).Value = 6479.4
$ws.Cells.Item(7, 11).Value = 3873.2727
$ws.Cells.Item(7, 12).Value = 6479.4
$ws.Cells.Item(7, 13).Value = -3761.2727
$ws.Cells.Item(7, 14).Value = -6703.4
$ws.Cells.Item(61, 8).Value = 2486.75
$ws.Cells.Item(61, 9).Value = 2766
$ws.Cells.Item(61, 11).Value = 2766
$ws.Cells.Item(61, 13).Value = -2564
$ws.Cells.Item(113, 8).Value = 2486.75
$ws.Cells.Item(113, 9).Value = 2766
$ws.Cells.Item(113, 11).Value = 2766
$ws.Cells.Item(113, 13).Value = -596
$ws.Cells.Item(122, 8).Value = 3359.25
$ws.Cells.Item(122, 9).Value = 3359.25
$ws.Cells.Item(122, 11).Value = 10077.75
$ws.Cells.Item(122, 13).Value = -7627.75
$ws.Cells.Item(126, 8).Value = 4687.6875
$ws.Cells.Item(126, 9).Value = 3873.2727
$ws.Cells.Item(126, 10).Value = 6479.4
$ws.Cells.Item(126, 11).Value = 11619.8181
$ws.Cells.Item(126, 12).Value = 19438.2
$ws.Cells.Item(126, 13).Value = -9149.8181
$ws.Cells.Item(126, 14).Value = -24378.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 1136.0303
$ws.Cells.Item(18, 9).Value = 1136.0303
$ws.Cells.Item(18, 11).Value = 1136.0303
$ws.Cells.Item(18, 13).Value = -963.0302999999999
$ws.Cells.Item(81, 8).Value = 4318.7144
$ws.Cells.Item(81, 9).Value = 4572.4
$ws.Cells.Item(81, 10).Value = 3684.5
$ws.Cells.Item(81, 11).Value = 9144.799999999999
$ws.Cells.Item(81, 12).Value = 7369
$ws.Cells.Item(81, 13).Value = -8083.799999999999
$ws.Cells.Item(81, 14).Value = -9491
$ws.Cells.Item(84, 8).Value = 4318.7144
$ws.Cells.Item(84, 9).Value = 4572.4
$ws.Cells.Item(84, 10).Value = 3684.5
$ws.Cells.Item(84, 11).Value = 45724
$ws.Cells.Item(84, 12).Value = 36845
$ws.Cells.Item(84, 13).Value = -40420
$ws.Cells.Item(84, 14).Value = -47453
$ws.Cells.Item(96, 8).Value = 4398
$ws.Cells.Item(96, 9).Value = 3833
$ws.Cells.Item(96, 10).Value = 4963
$ws.Cells.Item(96, 11).Value = 3833
$ws.Cells.Item(96, 12).Value = 4963
$ws.Cells.Item(96, 13).Value = -2460
$ws.Cells.Item(96, 14).Value = -7709
$ws.Cells.Item(122, 8).Value = 6606
$ws.Cells.Item(122, 9).Value = 6606
$ws.Cells.Item(122, 11).Value = 19818
$ws.Cells.Item(122, 13).Value = -17368
$ws.Cells.Item(126, 8).Value = 1198.3529
$ws.Cells.Item(126, 9).Value = 1119.2667
$ws.Cells.Item(126, 10).Value = 1791.5
$ws.Cells.Item(126, 11).Value = 3357.800099999999
$ws.Cells.Item(126, 12).Value = 5374.5
$ws.Cells.Item(126, 13).Value = -887.8000999999995
$ws.Cells.Item(126, 14).Value = -10314.5
